# Update cryptocurrency price/volume data per Mon Jan  8 11:59:26 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.456.84'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '2.233.16'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.85%  '
$ws.Range('D5').Value = "'298.26"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.64%  '
$ws.Range('D6').Value = "'91.36"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.71%  '
$ws.Range('D7').Value = "'0.561"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.77%  '
$ws.Range('D8').Value = "'0.998"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').Value = "'0.498"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -4.51%  '
$ws.Range('D10').Value = "'33.55"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.91%  '
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('D12').Value = "'7.05"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('D13').Value = "'0.104"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '2.568.92'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').Value = '2.232.60'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').Value = "'13.41"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = "'0.782"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -6.13%  '
$ws.Range('D18').Value = '44.222.61'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = "'12.24"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.78%  '
$ws.Range('E20').Value = '  -4.45%  '
$ws.Range('D21').Value = "'6.03"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('D22').Value = "'64.51"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').Value = "'237.87"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = "'2.83"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.20%  '
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = "'1.85"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.09%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = "'2.27"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = "'38.94"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('D29').Value = "'9.42"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').Value = "'19.36"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.83%  '
$ws.Range('D31').Value = "'152.39"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  -7.67%  '
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('D34').Value = "'2.51"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.29%  '
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('D36').Value = "'0.105"
$ws.Range('D36').ClearFormats()
$ws.Range('D37').Value = "'2.87"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.14%  '
$ws.Range('E38').Value = '  -7.36%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('E40').Value = '  -5.97%  '
$ws.Range('D41').Value = "'3.62"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.07%  '
$ws.Range('D42').Value = "'13.57"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -9.54%  '
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').Value = '1.796.88'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('E45').Value = '  +8.62%  '
$ws.Range('D46').Value = "'0.185"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.34%  '
$ws.Range('D47').Value = "'69.15"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('D48').Value = "'95.26"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.52%  '
$ws.Range('D49').Value = "'74.12"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.75%  '
$ws.Range('D50').Value = "'4.64"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.45%  '
$ws.Range('D51').Value = "'7.78"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.79%  '
